$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.278141666666667
$ws.Range("H2").Value = 12.834425
$ws.Range("I2").Value = 0.9663225094340192
$ws.Range("J2").Value = 0.9663225094340191
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.695610666666666
$ws.Range("N2").Value = 14.086832
$ws.Range("O2").Value = 0.1802066564018305
$ws.Range("P2").Value = 0.1802066564018305
$ws.Range("Q2").Value = 20.08848764351111
$ws.Range("R2").Value = 180.7963887916
$ws.Range("S2").Value = 0.1741377484309309
$ws.Range("T2").Value = 0.1741377484309309

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.278141666666667
$ws.Range("H3").Value = 12.834425
$ws.Range("I3").Value = 0.9663225094340192
$ws.Range("J3").Value = 0.9663225094340191
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.51448033333333
$ws.Range("N3").Value = 46.543441
$ws.Range("O3").Value = 0.5954098039960916
$ws.Range("P3").Value = 0.5954098039960916
$ws.Range("Q3").Value = 66.37314475071389
$ws.Range("R3").Value = 597.358302756425
$ws.Range("S3").Value = 0.5753578959391207
$ws.Range("T3").Value = 0.5753578959391207

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.278141666666667
$ws.Range("H4").Value = 12.834425
$ws.Range("I4").Value = 0.9663225094340192
$ws.Range("J4").Value = 0.9663225094340191
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.846719333333333
$ws.Range("N4").Value = 17.540158
$ws.Range("O4").Value = 0.2243835396020779
$ws.Range("P4").Value = 0.2243835396020779
$ws.Range("Q4").Value = 25.01309359323889
$ws.Range("R4").Value = 225.11784233915
$ws.Range("S4").Value = 0.2168268650639676
$ws.Range("T4").Value = 0.2168268650639675

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1490983333333333
$ws.Range("H5").Value = 0.447295
$ws.Range("I5").Value = 0.03367749056598092
$ws.Range("J5").Value = 0.03367749056598091
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.695610666666666
$ws.Range("N5").Value = 14.086832
$ws.Range("O5").Value = 0.1802066564018305
$ws.Range("P5").Value = 0.1802066564018305
$ws.Range("Q5").Value = 0.7001077243822221
$ws.Range("R5").Value = 6.30096951944
$ws.Range("S5").Value = 0.006068907970899612
$ws.Range("T5").Value = 0.00606890797089961

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1490983333333333
$ws.Range("H6").Value = 0.447295
$ws.Range("I6").Value = 0.03367749056598092
$ws.Range("J6").Value = 0.03367749056598091
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 15.51448033333333
$ws.Range("N6").Value = 46.543441
$ws.Range("O6").Value = 0.5954098039960916
$ws.Range("P6").Value = 0.5954098039960916
$ws.Range("Q6").Value = 2.313183160232778
$ws.Range("R6").Value = 20.818648442095
$ws.Range("S6").Value = 0.02005190805697092
$ws.Range("T6").Value = 0.02005190805697092

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1490983333333333
$ws.Range("H7").Value = 0.447295
$ws.Range("I7").Value = 0.03367749056598092
$ws.Range("J7").Value = 0.03367749056598091
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.846719333333333
$ws.Range("N7").Value = 17.540158
$ws.Range("O7").Value = 0.2243835396020779
$ws.Range("P7").Value = 0.2243835396020779
$ws.Range("Q7").Value = 0.8717361080677777
$ws.Range("R7").Value = 7.845624972609999
$ws.Range("S7").Value = 0.007556674538110385
$ws.Range("T7").Value = 0.007556674538110382

